$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.612.74"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "1.869.29"
$ws.Range("E3").Value = "  +2.16%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.52"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4615"
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3876"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07868"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9750"
$ws.Range("E10").Value = "  +1.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.93"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "1.870.91"
$ws.Range("E12").Value = "  +4.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.986"
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.696"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06955"
$ws.Range("E15").Value = "  +2.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.13"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.005"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.80"
$ws.Range("E19").Value = "  +1.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").Value = "28.610.50"
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.272"
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.02"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.111"
$ws.Range("D25").Value = "2.085.14"
$ws.Range("E25").Value = "  +2.99%  "
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.761"
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.987"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.11"
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09337"
$ws.Range("E31").Value = "  +0.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9173"
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.335"
$ws.Range("E34").Value = "  +1.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.326"
$ws.Range("E35").Value = "  +1.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05793"
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.154"
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02090"
$ws.Range("E38").Value = "  -2.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.735"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5627"
$ws.Range("E40").Value = "  +0.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1784"
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.784"
$ws.Range("E42").Value = "  -0.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07178"
$ws.Range("E43").Value = "  +2.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.77"
$ws.Range("E44").Value = "  +1.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5300"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.165"
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.142"
$ws.Range("E47").Value = "  +3.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.831"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.92"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.407"
$ws.Range("E50").Value = "  +3.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.003"
$ws.Range("E51").Value = "  +0.20%  "
